$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D23").NumberFormat = "@"
$ws.Range("D25:D50").NumberFormat = "@"
$ws.Range("D2").Value = "27.073.48"
$ws.Range("E2").Value = "  -2.77%  "
$ws.Range("D3").Value = "1.715.77"
$ws.Range("E3").Value = "  -2.88%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "310.92"
$ws.Range("E5").Value = "  -5.04%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("D7").Value = "0.4592"
$ws.Range("E7").Value = "  +2.41%  "
$ws.Range("D8").Value = "0.3416"
$ws.Range("E8").Value = "  -3.70%  "
$ws.Range("D9").Value = "41.90"
$ws.Range("E9").Value = "  -0.47%  "
$ws.Range("D10").Value = "0.07241"
$ws.Range("E10").Value = "  -2.72%  "
$ws.Range("D11").Value = "1.039"
$ws.Range("E11").Value = "  -4.94%  "
$ws.Range("D12").Value = "1.004"
$ws.Range("E12").Value = "  +0.32%  "
$ws.Range("D13").Value = "19.69"
$ws.Range("E13").Value = "  -5.35%  "
$ws.Range("D14").Value = "5.819"
$ws.Range("E14").Value = "  -3.29%  "
$ws.Range("D15").Value = "1.722.38"
$ws.Range("E15").Value = "  -2.38%  "
$ws.Range("D16").Value = "6.850"
$ws.Range("E16").Value = "  -4.70%  "
$ws.Range("D17").Value = "88.31"
$ws.Range("E17").Value = "  -5.04%  "
$ws.Range("D18").Value = "0.00001037"
$ws.Range("E18").Value = "  -1.65%  "
$ws.Range("D19").Value = "0.06321"
$ws.Range("E19").Value = "  -1.67%  "
$ws.Range("D20").Value = "1.002"
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("D21").Value = "16.46"
$ws.Range("E21").Value = "  -4.03%  "
$ws.Range("D22").Value = "5.595"
$ws.Range("E22").Value = "  -3.02%  "
$ws.Range("D23").Value = "27.154.95"
$ws.Range("E23").Value = "  -2.63%  "
$ws.Range("E24").Value = "  -3.87%  "
$ws.Range("D25").Value = "2.128"
$ws.Range("E25").Value = "  +0.80%  "
$ws.Range("D26").Value = "154.64"
$ws.Range("E26").Value = "  -4.82%  "
$ws.Range("D27").Value = "19.30"
$ws.Range("E27").Value = "  -4.44%  "
$ws.Range("D28").Value = "1.917.98"
$ws.Range("E28").Value = "  -2.55%  "
$ws.Range("D29").Value = "2.129"
$ws.Range("E29").Value = "  -1.31%  "
$ws.Range("D30").Value = "119.80"
$ws.Range("E30").Value = "  -4.19%  "
$ws.Range("D31").Value = "1.021"
$ws.Range("E31").Value = "  -6.40%  "
$ws.Range("D32").Value = "0.09073"
$ws.Range("E32").Value = "  -0.62%  "
$ws.Range("D33").Value = "3.595"
$ws.Range("E33").Value = "  -1.62%  "
$ws.Range("D34").Value = "5.327"
$ws.Range("E34").Value = "  -4.15%  "
$ws.Range("D35").Value = "0.02190"
$ws.Range("E35").Value = "  -4.26%  "
$ws.Range("D36").Value = "0.05831"
$ws.Range("E36").Value = "  -4.35%  "
$ws.Range("D37").Value = "11.01"
$ws.Range("E37").Value = "  -6.92%  "
$ws.Range("D38").Value = "0.1988"
$ws.Range("E38").Value = "  -4.95%  "
$ws.Range("D39").Value = "4.713"
$ws.Range("E39").Value = "  -4.94%  "
$ws.Range("D40").Value = "1.401"
$ws.Range("E40").Value = "  +0.64%  "
$ws.Range("D41").Value = "0.5890"
$ws.Range("E41").Value = "  -6.15%  "
$ws.Range("D42").Value = "1.127"
$ws.Range("E42").Value = "  -4.70%  "
$ws.Range("D43").Value = "7.444"
$ws.Range("E43").Value = "  -5.97%  "
$ws.Range("D44").Value = "12.71"
$ws.Range("E44").Value = "  -3.82%  "
$ws.Range("D45").Value = "3.580"
$ws.Range("E45").Value = "  -4.19%  "
$ws.Range("D46").Value = "0.5615"
$ws.Range("E46").Value = "  -3.92%  "
$ws.Range("D47").Value = "118.68"
$ws.Range("E47").Value = "  -2.95%  "
$ws.Range("D48").Value = "1.858"
$ws.Range("E48").Value = "  -4.33%  "
$ws.Range("D49").Value = "0.06655"
$ws.Range("E49").Value = "  -3.65%  "
$ws.Range("D50").Value = "1.080"
$ws.Range("E50").Value = "  -4.72%  "
$ws.Range("E51").Value = "  +0.17%  "
